# Updates the "cryptos" price list worksheet with refreshed price/volume
# data pulled by the periodic GitHub Actions scrape. A couple of coins
# (TRON/WrappedEther at rows 16-17 and Stellar/ApeXProtocol at rows 47-48)
# swapped rank order, so their Coin/Link/Price/Volume cells are rewritten
# together; every other row just gets fresh Price (D) / Volume(1h) (E)
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells in column D are free-text (e.g. "66.017.56", "0.0000253")
# rather than numbers -- Excel's automatic type detection would otherwise
# coerce plain numeric-looking strings like "570.71" into real numbers
# (and group-separated ones like "66.017.56" would throw/garble). Forcing
# the cell to Text before the write keeps it a string, and resetting the
# style back to "Normal" afterwards avoids leaving a stray number-format
# override on the cell.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "66.017.56"
$ws.Range("E2").Value = "  -4.38%  "
Set-TextValue $ws.Range("D3") "3.522.79"
$ws.Range("E3").Value = "  -5.31%  "
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue $ws.Range("D5") "570.71"
$ws.Range("E5").Value = "  -6.65%  "
Set-TextValue $ws.Range("D6") "187.86"
$ws.Range("E6").Value = "  -1.30%  "
Set-TextValue $ws.Range("D7") "3.520.63"
$ws.Range("E7").Value = "  -5.26%  "
Set-TextValue $ws.Range("D8") "0.606"
$ws.Range("E8").Value = "  -4.96%  "
$ws.Range("E9").Value = "  +0.19%  "
Set-TextValue $ws.Range("D10") "0.659"
$ws.Range("E10").Value = "  -8.30%  "
$ws.Range("E11").Value = "  -10.79%  "
Set-TextValue $ws.Range("D12") "52.33"
$ws.Range("E12").Value = "  -10.30%  "
Set-TextValue $ws.Range("D13") "0.0000253"
$ws.Range("E13").Value = "  -12.68%  "
Set-TextValue $ws.Range("D14") "9.64"
$ws.Range("E14").Value = "  -9.16%  "
Set-TextValue $ws.Range("D15") "4.090.22"
$ws.Range("E15").Value = "  -5.27%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D16") "0.125"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D17") "3.524.64"
$ws.Range("E17").Value = "  -5.27%  "
Set-TextValue $ws.Range("D18") "18.07"
$ws.Range("E18").Value = "  -6.47%  "
Set-TextValue $ws.Range("D19") "65.818.67"
$ws.Range("E19").Value = "  -4.44%  "
Set-TextValue $ws.Range("D20") "11.97"
$ws.Range("E20").Value = "  -7.65%  "
Set-TextValue $ws.Range("D21") "1.04"
$ws.Range("E21").Value = "  -8.38%  "
Set-TextValue $ws.Range("D22") "388.45"
$ws.Range("E22").Value = "  -5.61%  "
Set-TextValue $ws.Range("D23") "4.23"
$ws.Range("E23").Value = "  -7.90%  "
Set-TextValue $ws.Range("D24") "84.41"
$ws.Range("E24").Value = "  -5.61%  "
Set-TextValue $ws.Range("D25") "10.88"
$ws.Range("E25").Value = "  -0.17%  "
Set-TextValue $ws.Range("D26") "2.86"
$ws.Range("E26").Value = "  -6.40%  "
Set-TextValue $ws.Range("D27") "12.17"
$ws.Range("E27").Value = "  -5.63%  "
$ws.Range("E28").Value = "  -0.25%  "
Set-TextValue $ws.Range("D29") "3.46"
$ws.Range("E29").Value = "  -9.01%  "
Set-TextValue $ws.Range("D30") "8.72"
$ws.Range("E30").Value = "  -9.73%  "
Set-TextValue $ws.Range("D31") "30.62"
$ws.Range("E31").Value = "  -7.57%  "
Set-TextValue $ws.Range("D32") "7.05"
$ws.Range("E32").Value = "  -5.74%  "
Set-TextValue $ws.Range("D33") "622.26"
$ws.Range("E33").Value = "  -0.58%  "
Set-TextValue $ws.Range("D34") "12.00"
Set-TextValue $ws.Range("D35") "63.06"
$ws.Range("E35").Value = "  -3.90%  "
Set-TextValue $ws.Range("D36") "0.111"
$ws.Range("E36").Value = "  -9.44%  "
Set-TextValue $ws.Range("D37") "40.95"
$ws.Range("E37").Value = "  -10.73%  "
$ws.Range("E38").Value = "  +0.13%  "
Set-TextValue $ws.Range("D39") "0.390"
$ws.Range("E39").Value = "  -5.80%  "
Set-TextValue $ws.Range("D40") "0.0₃0740"
$ws.Range("E40").Value = "  -9.64%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("E42").Value = "  -7.78%  "
Set-TextValue $ws.Range("D43") "2.949.63"
$ws.Range("E43").Value = "  +3.33%  "
Set-TextValue $ws.Range("D44") "2.76"
$ws.Range("E44").Value = "  -9.38%  "
Set-TextValue $ws.Range("D45") "2.44"
$ws.Range("E45").Value = "  -6.94%  "
Set-TextValue $ws.Range("D46") "0.0399"
$ws.Range("E46").Value = "  -10.44%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D47") "0.129"
$ws.Range("E47").Value = "  -8.04%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D48") "3.05"
$ws.Range("E48").Value = "  -1.31%  "
Set-TextValue $ws.Range("D49") "138.21"
$ws.Range("E49").Value = "  -3.83%  "
Set-TextValue $ws.Range("D50") "8.32"
$ws.Range("E50").Value = "  -8.63%  "
$ws.Range("E51").Value = "  -9.49%  "

